$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 268 ("さあおやすみ" post) entirely; this shifts all subsequent
# rows (269-396) up by one, matching the target state (A1:C395).
$ws.Rows.Item(268).Delete()
